# "added GraphPanel to MainGUI"
# Update the Sprint Backlog sheet: change status of two existing items and
# add two new Sprint Backlog rows (Edge Format / Vertex Format) tied to the
# newly added Graph Panel GUI work.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint Backlog")

# Row 2 ("Design Menue GUI Prototype") is now finished.
$ws.Range("K2").Value2 = "done"

# Row 6 ("Information Panel") has moved from "open" to "in Progess".
$ws.Range("K6").Value2 = "in Progess"

# New row 12: Edge Format story, owned by menzs2.
$ws.Range("D12").Value2 = "Edge Format"
$ws.Range("F12").Value2 = "menzs2"

# New row 13: Vertex Format story, owned by menzs2.
$ws.Range("D13").Value2 = "Vertex Format"
$ws.Range("F13").Value2 = "menzs2"

# Set statuses for the two new rows afterward (correct spelling "in Progress").
$ws.Range("K12").Value2 = "in Progress"
$ws.Range("K13").Value2 = "in Progress"

# Match the active selection recorded in the saved workbook.
$ws.Range("J12").Select()
